# Weekly update: insert a new week's price record for
# "Bruselas (repollito)" at Vega Central Mapocho de Santiago, pushing the
# existing history down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 34 (shifts rows 34:46 -> 35:47,
# carrying their formatting/styles with them, same as Excel's
# Rows.Insert()).
$ws.Rows("34:34").Insert()

# Populate the newly inserted row 34 with this week's data.
$ws.Range("A34").Value = 9
$ws.Range("B34").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value = 44755
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 100112035
$ws.Range("G34").Value = "Bruselas (repollito)"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 43
$ws.Range("K34").Value = 18000
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = 19023
$ws.Range("N34").Value = "`$/malla 15 kilos"
$ws.Range("O34").Value = "Hijuelas"
$ws.Range("P34").Value = 1268
$ws.Range("Q34").Value = 15
$ws.Range("R34").Value = "Hortaliza"
